# Sabbineni_LabExam03Grading.xlsx - "Changes to driver class 65-80"
#
# The CustomerMappingDriver Class section (rows 27-31) is regraded:
#  - addProduct() method (row 29) score raised from 8 to 9 points, and its
#    grading comment is replaced to call out the missing Customer-object
#    initialization / hyphen-splitting issue.
#  - findNoOfCustomers() method (row 30) keeps its 0 points but gets a new
#    comment explaining the incorrect output.
#
# All of the section subtotal / grand total formulas (E31, E38, D38, ...)
# recompute automatically from the new E29 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: the grading-comments column picks up two brand new shared
# strings. F30's text was the first of the pair introduced, F29's the
# second, so write F30 before F29 to match that allocation order.
$ws.Range("F30").Value = "(-4) for incorrect output for all methods."
$ws.Range("F29").Value = "(-2) for not intializing Customer object, (-5) for not reading tehe products and their brands and splitting them by hyphen and adding them to the inventory"

# addProduct() method score: 8 -> 9
$ws.Range("E29").Value = 9

# Reflect the author's final selection/scroll position on the sheet.
$ws.Range("F30").Select() | Out-Null
